$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target row plus the new values for any of the B/C/D/E
# columns that change for that row (columns not present are left untouched).
$updates = @(
    @{Row=2; D='28.044.45'; E='  +1.83%  '}
    @{Row=3; D='1.864.12'; E='  +1.01%  '}
    @{Row=4; D='1.001'; E='  -0.42%  '}
    @{Row=5; D='335.96'; E='  +0.60%  '}
    @{Row=6; E='  -0.38%  '}
    @{Row=7; D='0.4703'; E='  +1.49%  '}
    @{Row=8; D='0.3904'; E='  +1.49%  '}
    @{Row=9; D='46.77'; E='  +1.79%  '}
    @{Row=10; D='0.07964'; E='  +1.02%  '}
    @{Row=11; D='0.9851'; E='  -1.14%  '}
    @{Row=12; D='21.55'; E='  +0.44%  '}
    @{Row=13; D='5.950'; E='  -0.07%  '}
    @{Row=14; D='1.847.92'; E='  +0.20%  '}
    @{Row=15; D='7.206'; E='  +1.10%  '}
    @{Row=16; D='91.51'; E='  +3.47%  '}
    @{Row=17; E='  -0.45%  '}
    @{Row=18; D='0.00001044'; E='  +0.96%  '}
    @{Row=19; D='0.06597'; E='  -1.10%  '}
    @{Row=20; D='17.57'; E='  +2.78%  '}
    @{Row=21; D='1.001'; E='  -0.41%  '}
    @{Row=22; D='28.045.69'; E='  +1.84%  '}
    @{Row=23; D='5.409'; E='  +0.40%  '}
    @{Row=24; D='10.99'; E='  +0.94%  '}
    @{Row=25; D='2.288'; E='  -0.98%  '}
    @{Row=26; B='WrappedliquidstakedEther2.0'; C='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D='2.071.87'; E='  +0.31%  '}
    @{Row=27; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='159.37'; E='  +0.48%  '}
    @{Row=28; D='19.55'; E='  +0.40%  '}
    @{Row=29; D='2.111'; E='  -0.12%  '}
    @{Row=30; D='5.494'; E='  +1.69%  '}
    @{Row=31; D='119.31'; E='  -0.43%  '}
    @{Row=32; D='0.9643'; E='  -1.22%  '}
    @{Row=33; D='0.09493'; E='  +1.04%  '}
    @{Row=34; D='3.578'; E='  -0.39%  '}
    @{Row=35; D='5.315'; E='  +0.28%  '}
    @{Row=36; D='1.351'; E='  +0.83%  '}
    @{Row=37; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.02261'; E='  +1.49%  '}
    @{Row=38; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.06100'; E='  +0.95%  '}
    @{Row=39; D='8.318'; E='  +0.41%  '}
    @{Row=40; D='1.164'; E='  -1.18%  '}
    @{Row=41; D='1.000'; E='  -0.29%  '}
    @{Row=42; D='0.5944'; E='  +0.83%  '}
    @{Row=43; D='0.1871'; E='  +0.54%  '}
    @{Row=44; D='10.23'; E='  -1.07%  '}
    @{Row=45; E='  +3.72%  '}
    @{Row=46; D='0.5571'; E='  -0.25%  '}
    @{Row=47; E='  -0.03%  '}
    @{Row=48; D='1.959'; E='  +2.65%  '}
    @{Row=49; D='0.06878'; E='  +2.71%  '}
    @{Row=50; D='111.61'; E='  +0.65%  '}
    @{Row=51; E='  -32.75%  '}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("B")) { $ws.Range("B$r").Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Range("C$r").Value = $u.C }
    if ($u.ContainsKey("D")) {
        # Price column cells are stored as plain text (e.g. "21.55" or
        # "27.992.20"). Force text entry via NumberFormat so Excel does not
        # reinterpret the string as a number, then restore the original
        # (General) cell style so only the content changes.
        $origStyle = $ws.Range("D$r").Style
        $ws.Range("D$r").NumberFormat = "@"
        $ws.Range("D$r").Value = $u.D
        $ws.Range("D$r").Style = $origStyle
    }
    if ($u.ContainsKey("E")) { $ws.Range("E$r").Value = $u.E }
}

Write-Output "Updated $($updates.Count) rows"